$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values in this sheet are stored as literal text (e.g. "314.80", "3.17%"),
# not numbers. Force text entry (NumberFormat "@") so Excel does not silently
# convert the numeric-looking strings / percents into real numbers, then reset
# the style back to Normal so no stray number-format style sticks to the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "314.54"
Set-TextValue $ws.Range("E2") "3.03%"
Set-TextValue $ws.Range("D3") "39.46"
Set-TextValue $ws.Range("E3") "2.60%"
Set-TextValue $ws.Range("D4") "5.112"
Set-TextValue $ws.Range("E4") "0.40%"
Set-TextValue $ws.Range("D5") "0.08192"
Set-TextValue $ws.Range("E5") "1.73%"
Set-TextValue $ws.Range("D6") "1.965"
Set-TextValue $ws.Range("E6") "1.67%"
Set-TextValue $ws.Range("D7") "8.245"
Set-TextValue $ws.Range("E7") "3.77%"
Set-TextValue $ws.Range("D8") "0.9324"
Set-TextValue $ws.Range("E8") "0.34%"
Set-TextValue $ws.Range("D9") "0.1412"
Set-TextValue $ws.Range("E9") "-1.68%"
Set-TextValue $ws.Range("D10") "0.1977"
Set-TextValue $ws.Range("E10") "2.87%"
Set-TextValue $ws.Range("D11") "0.09146"
Set-TextValue $ws.Range("E11") "1.67%"
Set-TextValue $ws.Range("D12") "0.03521"
Set-TextValue $ws.Range("E12") "0.39%"
Set-TextValue $ws.Range("E13") "0.67%"
Set-TextValue $ws.Range("D14") "0.001400"
Set-TextValue $ws.Range("E14") "0.24%"
Set-TextValue $ws.Range("D15") "0.006117"
Set-TextValue $ws.Range("E15") "1.75%"
Set-TextValue $ws.Range("D16") "3.663"
Set-TextValue $ws.Range("E16") "-1.63%"
Set-TextValue $ws.Range("D17") "4.275"
Set-TextValue $ws.Range("E17") "1.95%"
Set-TextValue $ws.Range("D18") "3.322"
Set-TextValue $ws.Range("E18") "-2.95%"
Set-TextValue $ws.Range("E19") "0.01%"
Set-TextValue $ws.Range("E20") "-2.37%"
Set-TextValue $ws.Range("D21") "4.867"
Set-TextValue $ws.Range("E21") "0.65%"
Set-TextValue $ws.Range("E22") "1.62%"
Set-TextValue $ws.Range("D23") "0.04319"
Set-TextValue $ws.Range("E23") "-0.93%"
Set-TextValue $ws.Range("D24") "0.001222"
Set-TextValue $ws.Range("E24") "-0.35%"
Set-TextValue $ws.Range("D25") "0.004789"
Set-TextValue $ws.Range("E25") "16.19%"
Set-TextValue $ws.Range("D27") "0.0003994"
Set-TextValue $ws.Range("E27") "-10.19%"
Set-TextValue $ws.Range("D39") "0.02249"
Set-TextValue $ws.Range("E39") "8.29%"
Set-TextValue $ws.Range("D40") "0.05287"
Set-TextValue $ws.Range("E40") "5.07%"
Set-TextValue $ws.Range("D41") "0.007611"
Set-TextValue $ws.Range("E41") "1.81%"
Set-TextValue $ws.Range("D42") "0.009781"
Set-TextValue $ws.Range("E42") "-3.45%"
Set-TextValue $ws.Range("D43") "0.1379"
Set-TextValue $ws.Range("E43") "2.20%"
Set-TextValue $ws.Range("D44") "0.002137"
Set-TextValue $ws.Range("E44") "-0.32%"
Set-TextValue $ws.Range("D45") "0.009801"
Set-TextValue $ws.Range("E45") "11.33%"
Set-TextValue $ws.Range("D46") "0.00006371"
Set-TextValue $ws.Range("E46") "2.95%"
Set-TextValue $ws.Range("E48") "-7.43%"
Set-TextValue $ws.Range("E49") "-25.18%"
